$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.541.91"
$ws.Range("E2").Value = "  +1.90%  "
$ws.Range("D3").Value = "2.152.99"
$ws.Range("E3").Value = "  +2.25%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.59"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("E6").Value = "  +0.57%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "62.53"
$ws.Range("E7").Value = "  +0.91%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.389"
$ws.Range("E9").Value = "  -0.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0839"
$ws.Range("E10").Value = "  -0.45%  "
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.77"
$ws.Range("E12").Value = "  -1.43%  "
$ws.Range("D13").Value = "2.473.78"
$ws.Range("E13").Value = "  +2.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.63"
$ws.Range("E14").Value = "  -1.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.803"
$ws.Range("E15").Value = "  +0.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.45"
$ws.Range("E16").Value = "  -0.72%  "
$ws.Range("D17").Value = "2.142.16"
$ws.Range("E17").Value = "  +3.73%  "
$ws.Range("D18").Value = "39.524.01"
$ws.Range("E18").Value = "  +1.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.51"
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.04"
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("E21").Value = "  +0.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.03"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.34"
$ws.Range("E24").Value = "  +1.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.32"
$ws.Range("E25").Value = "  -1.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.50"
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E27").Value = "  -2.89%  "
$ws.Range("E28").Value = "  +1.18%  "
$ws.Range("E29").Value = "  +0.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.57"
$ws.Range("E30").Value = "  +1.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.68"
$ws.Range("E31").Value = "  +4.64%  "
$ws.Range("E32").Value = "  +0.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.56"
$ws.Range("E33").Value = "  -0.30%  "
$ws.Range("E34").Value = "  -2.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.95"
$ws.Range("E35").Value = "  -2.78%  "
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("E37").Value = "  +7.72%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.38"
$ws.Range("E38").Value = "  +0.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.91"
$ws.Range("E40").Value = "  +18.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "102.53"
$ws.Range("E41").Value = "  +0.45%  "
$ws.Range("E42").Value = "  -1.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.59"
$ws.Range("E43").Value = "  -2.04%  "
$ws.Range("D44").Value = "1.512.24"
$ws.Range("E44").Value = "  -1.09%  "
$ws.Range("E45").Value = "  -0.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.83"
$ws.Range("E46").Value = "  +0.75%  "
$ws.Range("E47").Value = "  +0.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0917"
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("E50").Value = "  +0.76%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.358.24"
$ws.Range("E51").Value = "  +2.41%  "
